# Fix header labels on the existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet as the last sheet in the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Reuse the same header style (bold, centered, bordered) already used on the
# other two sheets, instead of re-deriving new style/font records.
$wsWeekly.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows
$ws3.Cells.Item(2,1).Value = 44941.99999999999
$ws3.Cells.Item(2,2).Value = 110
$ws3.Cells.Item(2,3).Value = 48.18827133841597
$ws3.Cells.Item(2,4).Value = 176.2879155029086
$ws3.Cells.Item(3,1).Value = 44948.99999999999
$ws3.Cells.Item(3,2).Value = 107
$ws3.Cells.Item(3,3).Value = 39.98398044203599
$ws3.Cells.Item(3,4).Value = 173.1723436414674
$ws3.Cells.Item(4,1).Value = 45109.99999999999
$ws3.Cells.Item(4,2).Value = 18
$ws3.Cells.Item(4,3).Value = -49.95771586214815
$ws3.Cells.Item(4,4).Value = 82.9868003773275
$ws3.Cells.Item(5,1).Value = 45116.99999999999
$ws3.Cells.Item(5,2).Value = 15
$ws3.Cells.Item(5,3).Value = -49.02975269469693
$ws3.Cells.Item(5,4).Value = 83.12680227582335
$ws3.Cells.Item(6,1).Value = 45123.99999999999
$ws3.Cells.Item(6,2).Value = 11
$ws3.Cells.Item(6,3).Value = -55.88334866103106
$ws3.Cells.Item(6,4).Value = 75.29808248807737
$ws3.Cells.Item(7,1).Value = 45130.99999999999
$ws3.Cells.Item(7,2).Value = 7
$ws3.Cells.Item(7,3).Value = -61.5933558867499
$ws3.Cells.Item(7,4).Value = 72.46258431452709
$ws3.Cells.Item(8,1).Value = 45137.99999999999
$ws3.Cells.Item(8,2).Value = 3
$ws3.Cells.Item(8,3).Value = -63.31571478883333
$ws3.Cells.Item(8,4).Value = 67.79087993288989
$ws3.Cells.Item(9,1).Value = 45144.99999999999
$ws3.Cells.Item(9,2).Value = 0
$ws3.Cells.Item(9,3).Value = -69.3811259720115
$ws3.Cells.Item(9,4).Value = 66.30880630135319
$ws3.Cells.Item(10,1).Value = 45151.99999999999
$ws3.Cells.Item(10,2).Value = 0
$ws3.Cells.Item(10,3).Value = -72.19636091253452
$ws3.Cells.Item(10,4).Value = 61.62718605499659
$ws3.Cells.Item(11,1).Value = 45158.99999999999
$ws3.Cells.Item(11,2).Value = 0
$ws3.Cells.Item(11,3).Value = -73.65387954735577
$ws3.Cells.Item(11,4).Value = 58.56190925668568
$ws3.Cells.Item(12,1).Value = 45165.99999999999
$ws3.Cells.Item(12,2).Value = 0
$ws3.Cells.Item(12,3).Value = -78.71891838069639
$ws3.Cells.Item(12,4).Value = 49.65431352457735
$ws3.Cells.Item(13,1).Value = 45172.99999999999
$ws3.Cells.Item(13,2).Value = 0
$ws3.Cells.Item(13,3).Value = -82.37879357075548
$ws3.Cells.Item(13,4).Value = 45.03974309085527
$ws3.Cells.Item(14,1).Value = 45179.99999999999
$ws3.Cells.Item(14,2).Value = 0
$ws3.Cells.Item(14,3).Value = -85.13591131300173
$ws3.Cells.Item(14,4).Value = 48.63448525600768
$ws3.Cells.Item(15,1).Value = 45186.99999999999
$ws3.Cells.Item(15,2).Value = 0
$ws3.Cells.Item(15,3).Value = -86.39147117286039
$ws3.Cells.Item(15,4).Value = 40.57049785338592
$ws3.Cells.Item(16,1).Value = 45193.99999999999
$ws3.Cells.Item(16,2).Value = 0
$ws3.Cells.Item(16,3).Value = -97.52211458493579
$ws3.Cells.Item(16,4).Value = 33.24509844909836

# Reuse the existing date-time display style on the "ds" column, matching
# the "Order Week"/"Order Month" date columns on the other two sheets.
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
